# Update "想去人数" (interest count) values in column F across sheets
# 展览(1), 演出(2), 本地生活(3), 全部类型(4)
# to reflect the refreshed scrape at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 832
$ws.Range("F4").Value = 429
$ws.Range("F5").Value = 882
$ws.Range("F6").Value = 493
$ws.Range("F7").Value = 7371
$ws.Range("F8").Value = 140
$ws.Range("F10").Value = 1929
$ws.Range("F11").Value = 5377
$ws.Range("F15").Value = 7434
$ws.Range("F16").Value = 8727
$ws.Range("F19").Value = 864
$ws.Range("F20").Value = 4369
$ws.Range("F22").Value = 192
$ws.Range("F26").Value = 1179
$ws.Range("F28").Value = 1620
$ws.Range("F29").Value = 692
$ws.Range("F31").Value = 1850
$ws.Range("F32").Value = 318
$ws.Range("F33").Value = 2239
$ws.Range("F34").Value = 313
$ws.Range("F36").Value = 1418
$ws.Range("F39").Value = 782
$ws.Range("F40").Value = 384
$ws.Range("F41").Value = 4022
$ws.Range("F44").Value = 406
$ws.Range("F48").Value = 160
$ws.Range("F49").Value = 4053
$ws = $wb.Worksheets.Item(2)
$ws.Range("F14").Value = 45
$ws.Range("F25").Value = 59
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5104
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 5104
$ws.Range("F5").Value = 832
$ws.Range("F6").Value = 882
$ws.Range("F7").Value = 493
$ws.Range("F11").Value = 5377
$ws.Range("F13").Value = 7434
$ws.Range("F17").Value = 864
$ws.Range("F18").Value = 4369
$ws.Range("F20").Value = 192
$ws.Range("F25").Value = 1179
$ws.Range("F27").Value = 1620
$ws.Range("F28").Value = 1850
$ws.Range("F29").Value = 318
$ws.Range("F30").Value = 2239
$ws.Range("F37").Value = 782
$ws.Range("F39").Value = 59
$ws.Range("F40").Value = 384
$ws.Range("F41").Value = 4022
$ws.Range("F45").Value = 406
$ws.Range("F48").Value = 160
$ws.Range("F49").Value = 4053
